# Update "想去人数" (interest count) and "最低票价" (lowest price) figures
# for the 合肥·第九届环形宇宙动漫游戏嘉年华 event, on both the "展览" sheet
# and the consolidated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("G2").Value = 65
$wsExhibit.Range("F7").Value = 1602

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("G2").Value = 65
$wsAll.Range("F11").Value = 1602
